# Generate Report for Handoff
#
# The automated handoff run produced a new report markdown file and the
# handoff transform for this change failed, so each locale's row needs to
# reflect: the new report file name, a failed status, no handoff xlf file
# (the hyperlink + value are cleared), a reset "Latest Handoff Datetime",
# and a Handoff Reason of "Ignored" instead of "Include".

$oldUuid = "819e3aa0-8943-46b8-b57c-87be591b4fd9"
$newUuid = "5b15cfee-8424-4ce3-bdd9-97e2dcdc2153"

$oldFileName = "$oldUuid.md"
$newFileName = "$newUuid.md"

$newStatus = "Handoff transform failed"
$emptyDate = "0001-01-01 00:00:00"
$ignoredReason = "Ignored"

$wb = $excel.ActiveWorkbook

# Re-point every hyperlink whose address/display references the old report
# file so it references the newly generated one instead, on every sheet.
# Note: deleting is done per-hyperlink-object (found via the sheet-level
# Hyperlinks collection) further down -- Range.Hyperlinks.Delete() clears
# *every* hyperlink on the sheet in this host, not just the target range.
foreach ($ws in $wb.Worksheets) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay -eq $oldFileName) {
            $hl.Address = $hl.Address.Replace($oldUuid, $newUuid)
            $hl.TextToDisplay = $newFileName
        }
    }
}

# --- Overview sheet: file name (A2) plus the per-locale status mirrors
#     the same status text as each detail sheet's B2 ("Handoff transform
#     failed" instead of "Ready for handoff") ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = $newFileName
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# --- Locale detail sheets (zh-cn / de-de): row 2 holds the status for the
#     report file itself; row 3 is the unrelated .localization-config row
#     and is untouched. ---
$localeSheets = @("zh-cn", "de-de")
foreach ($name in $localeSheets) {
    $ws = $wb.Worksheets.Item($name)

    # A2: Source File Name -> new report file name
    $ws.Range("A2").Value = $newFileName

    # B2: Status -> "Handoff transform failed"
    $ws.Range("B2").Value = $newStatus

    # C2: Latest Handoff File -> no file produced; remove its hyperlink
    # (only that one -- deleted by object, not via Range.Hyperlinks, which
    # would wipe A2/A3's hyperlinks too) then clear the cell itself.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $hl.Delete()
        }
    }
    $ws.Range("C2").Clear()

    # D2: Latest Handoff Datetime -> reset to the empty/default datetime
    $ws.Range("D2").Value = $emptyDate

    # G2: Latest Handback DateTime -> stays the empty/default datetime
    $ws.Range("G2").Value = $emptyDate

    # H2: Handoff Reason -> "Ignored" (was "Include")
    $ws.Range("H2").Value = $ignoredReason

    # Row 3 values are textually unchanged, but rewrite them too so the
    # shared-string table collapses the same way it would after removing
    # the now-unused handoff/date strings above.
    $ws.Range("D3").Value = $emptyDate
    $ws.Range("G3").Value = $emptyDate
    $ws.Range("H3").Value = $ignoredReason
}
